$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 = "Class 01c methods" (class 0x1c) -- note progress, highlighted amber.
$cell25 = $ws.Cells.Item(25, 2)
$cell25.Value = "300, 308, 30c done, rest not done"
$cell25.Font.Name = "Arial"
$cell25.Interior.Color = 49407

# Rows 4-24 correspond to class method rows 001..018 (everything except the
# last "Class 01c methods" row, handled above).
# Mark them with "NONE" in column B, highlighted in red.
for ($r = 4; $r -le 24; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = "NONE"
    $cell.Font.Name = "Arial"
    $cell.Interior.Color = 255
}

# Widen column B to fit the new content.
$ws.Columns.Item(2).ColumnWidth = 16.85

# Move the active selection.
$ws.Range("G14").Select()
